$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Simple value updates (rows 3, 6, 8) ---
$ws.Range("B3").Value2 = "0.1.7"
$ws.Range("B6").Value2 = "draft"
$ws.Range("B8").Value2 = "2024-11-22T12:33:30-06:00"

# --- Contact section rework ---
# Row 10 Contact value is replaced with the publisher/contact org text.
$ws.Range("B10").Value2 = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# Prep a new row 16 whose format (style) matches row 15, so the shift below
# keeps every cell on style "s=2" instead of Excel minting a brand-new xf.
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)

# Shift old rows 12-15 down to 13-16 (bottom-up, so we never clobber a value
# before we've read it).
$ws.Range("A16").Value2 = $ws.Range("A15").Value2
$ws.Range("B16").Value2 = $ws.Range("B15").Value2

$ws.Range("A15").Value2 = $ws.Range("A14").Value2
$ws.Range("B15").Value2 = $ws.Range("B14").Value2

$ws.Range("A14").Value2 = $ws.Range("A13").Value2
$ws.Range("B14").Value2 = $ws.Range("B13").Value2

$ws.Range("A13").Value2 = $ws.Range("A12").Value2
$ws.Range("B13").Value2 = $ws.Range("B12").Value2

# Row 11 (previously a duplicate "Contact" row) now holds the second contact.
$ws.Range("A11").Value2 = "Contact"
$ws.Range("B11").Value2 = "Bob Milius (bmilius@nmdp.org)"

# Row 12 becomes the new "Jurisdiction" property with a blank value.
$ws.Range("A12").Value2 = "Jurisdiction"
$ws.Range("B12").Value2 = ""

$wb.Save()
